# Fruta / hortaliza, semanal
# Insert a new weekly record at row 86 (Florida King, Provincia de Limarí),
# pushing the existing rows 86-105 down to 87-106.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 86:105 down by one to make room for the new record.
$ws.Rows.Item(86).EntireRow.Insert()

# Populate the newly inserted row 86 with the new weekly record.
$row = 86
$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 44510
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100103
$ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value = 100103004
$ws.Cells.Item($row, 10).Value = "Durazno"
$ws.Cells.Item($row, 11).Value = "Florida King"
$ws.Cells.Item($row, 12).Value = "Segunda"
$ws.Cells.Item($row, 13).Value = 220
$ws.Cells.Item($row, 14).Value = 15000
$ws.Cells.Item($row, 15).Value = 16000
$ws.Cells.Item($row, 16).Value = 15545
$ws.Cells.Item($row, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 19).Value = 1036
$ws.Cells.Item($row, 20).Value = 15
